# Update "想去人数" (want-to-go count) figures for several events.
# These updates apply to both the "展览" sheet and the aggregated
# "全部类型" sheet, which lists the same events at different rows.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 5856
$ws1.Range("F6").Value = 5179
$ws1.Range("F7").Value = 435
$ws1.Range("F12").Value = 34

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 5856
$ws4.Range("F6").Value = 5179
$ws4.Range("F7").Value = 435
$ws4.Range("F14").Value = 34
